$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hours for "Investigación y documentación app" (row 6) and "Creación app" (row 7)
$ws.Range("E6").Value = 16
$ws.Range("E7").Value = 34

# Move selection/active cell to E8 (no results screen added task row)
$ws.Range("E8").Select()
